$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark from its old location
#    (it sat right after "...As palavras desse dicionário ").
#    We will re-create it later, at the end of the investments
#    paragraph, once all text edits are done (so its offset is
#    correct for the final document).
# ---------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------
# 2) Rewrite the "Toda notícia ..." paragraph with the expanded
#    description (likes/dislikes, localidade, fonte, manchete,
#    link da imagem, highlight-selection paragraph, etc.)
# ---------------------------------------------------------------
$oldNews = "Toda notícia contém comentários, curtidas, link que redireciona para ela e um ID. Todo comentário tem uma data e hora, curtidas, quem e o que comentou, além das respostas que também são outros comentários."
$newNews = "Toda notícia contém comentários, curtidas e não curtidas, link que redireciona para ela e um ID, localidade, fonte, descrição, manchete, fonte, link para a imagem e, em alguns casos, o conteúdo. Para escolher qual notícia ficará em destaque, a quantidade de curtidas, de comentários e cliques aumentam o engajamento, quanto maior mais provável que seja destaque. Todo comentário tem uma data e hora, curtidas, quem e o que comentou, além das respostas que também são outros comentários."

$d.Content.Find.Execute($oldNews, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newNews, 2) | Out-Null

# ---------------------------------------------------------------
# 3) Rewrite the investments-simulation paragraph: drop the
#    "de fim (...) e o tipo (poupança, bolsa, tesouro etc)" part
#    and replace the old closing sentence with the new text about
#    "Toda simulação é de um investimento ...".
# ---------------------------------------------------------------
$oldSim = "data de início da simulação, de fim (que depende do tipo de investimento que o usuário escolheu) e o tipo (poupança, bolsa, tesouro etc) que deve armazenar a descrição, rendimento e o período que vai render."
$newSim = "data de início da simulação. Toda simulação é de um investimento, que tem um nome, uma descrição, grupo, período, rendimento e a quantidade de simulações daquele investimento para que se possa ter uma noção dos interesses dos usuários."

$d.Content.Find.Execute($oldSim, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newSim, 2) | Out-Null

# ---------------------------------------------------------------
# 4) Re-insert the "_GoBack" bookmark as a collapsed range right
#    at the end of the (now rewritten) investments paragraph -
#    i.e. right before its paragraph mark.
#
#    A range collapsed exactly on "paragraph.End - 1" cannot be
#    fed straight into Bookmarks.Add, so we insert a one-char
#    placeholder there, bookmark the range that contains it, and
#    then delete the placeholder again; the bookmark collapses
#    to the correct spot and survives the deletion.
# ---------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$endPos = $lastPara.Range.End - 1

$marker = $d.Range($endPos, $endPos)
$marker.InsertAfter("|")

$markerSpan = $d.Range($endPos, $endPos + 1)
$d.Bookmarks.Add("_GoBack", $markerSpan) | Out-Null

$markerSpan2 = $d.Range($endPos, $endPos + 1)
$markerSpan2.Text = ""
